# Realestate Update resale numbers 2025-01-05 11:44
# Append a new data row (row 9) to the CityResaleNum sheet with the
# day's resale figures.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row = 9

# Columns A and D hold text that looks numeric/date-like ("2025-01-05",
# "01"); force Text format first so COM doesn't coerce them into a date
# serial / plain number, then restore the default "Normal" style so the
# new row's cells stay unstyled like the rest of the data rows.
$ws.Range("A$row").NumberFormat = "@"
$ws.Range("D$row").NumberFormat = "@"

$ws.Cells.Item($row, 1).Value = "2025-01-05"
$ws.Cells.Item($row, 2).Value = "11:44:14"
$ws.Cells.Item($row, 3).Value = "Sunday"
$ws.Cells.Item($row, 4).Value = "01"

$ws.Range("A$row").Style = "Normal"
$ws.Range("D$row").Style = "Normal"

$ws.Cells.Item($row, 5).Value = 127672
$ws.Cells.Item($row, 6).Value = 143737
$ws.Cells.Item($row, 7).Value = 168493
$ws.Cells.Item($row, 8).Value = 158396
$ws.Cells.Item($row, 9).Value = -1
$ws.Cells.Item($row, 10).Value = 142221
$ws.Cells.Item($row, 11).Value = -1
$ws.Cells.Item($row, 12).Value = -1
$ws.Cells.Item($row, 13).Value = 192538
$ws.Cells.Item($row, 14).Value = 114874
$ws.Cells.Item($row, 15).Value = 45498
$ws.Cells.Item($row, 16).Value = 28311
$ws.Cells.Item($row, 17).Value = 63570
$ws.Cells.Item($row, 18).Value = -1
$ws.Cells.Item($row, 19).Value = 47987
$ws.Cells.Item($row, 20).Value = -1
